$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.163.14'
$ws.Range('E2').Value = '  +4.85%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.500.50'
$ws.Range('E3').Value = '  +3.08%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '323.86'
$ws.Range('E5').Value = '  +1.76%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.29'
$ws.Range('E6').Value = '  +1.11%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.521'
$ws.Range('E7').Value = '  +0.95%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.543'
$ws.Range('E9').Value = '  +2.17%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '37.11'
$ws.Range('E10').Value = '  +3.75%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0814'
$ws.Range('E11').Value = '  +1.28%  '
$ws.Range('E12').Value = '  +0.68%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.30'
$ws.Range('E13').Value = '  -0.99%  '
$ws.Range('E14').Value = '  +3.95%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.880.42'
$ws.Range('E15').Value = '  +2.70%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.547.51'
$ws.Range('E16').Value = '  +4.83%  '
$ws.Range('E17').Value = '  +1.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.043.69'
$ws.Range('E18').Value = '  +4.84%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.68'
$ws.Range('E19').Value = '  +2.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.62'
$ws.Range('E20').Value = '  +4.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0936'
$ws.Range('E21').Value = '  +1.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.61'
$ws.Range('E22').Value = '  +2.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '251.02'
$ws.Range('E24').Value = '  +2.98%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('E25').Value = '  +1.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.25'
$ws.Range('E26').Value = '  +3.46%  '
$ws.Range('E27').Value = '  -0.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.04'
$ws.Range('E28').Value = '  +4.83%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.19'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.14'
$ws.Range('E30').Value = '  +3.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.134'
$ws.Range('E31').Value = '  +3.73%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.52'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.67'
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('E34').Value = '  +1.58%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0779'
$ws.Range('E35').Value = '  +1.95%  '
$ws.Range('E36').Value = '  +0.07%  '
$ws.Range('E37').Value = '  +1.66%  '
$ws.Range('E38').Value = '  +1.72%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.97'
$ws.Range('E39').Value = '  +3.16%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.111'
$ws.Range('E40').Value = '  +1.42%  '
$ws.Range('B41').Value = 'Monero'
$ws.Range('C41').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '122.26'
$ws.Range('E41').Value = '  -2.77%  '
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.45'
$ws.Range('E43').Value = '  +1.35%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0296'
$ws.Range('E44').Value = '  +1.93%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.964.78'
$ws.Range('E45').Value = '  +0.93%  '
$ws.Range('E46').Value = '  +1.30%  '
$ws.Range('E47').Value = '  +0.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.80'
$ws.Range('E48').Value = '  +0.95%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.11'
$ws.Range('E49').Value = '  -1.27%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.42'
$ws.Range('E50').Value = '  +15.83%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '78.87'
$ws.Range('E51').Value = '  +4.09%  '
